$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 467 ("「これ、私の眼鏡じゃありません」" post) and shift
# everything below it up by one row, reducing the used range from
# A1:C589 to A1:C588.
$ws.Rows.Item(467).Delete()
